# Actualización automática 2025-08-19 15:45:09
#
# Registra una venta negativa (nota de crédito) de -81.41 para el cliente
# "ZAMBRANO ANGELA MARIA" en la categoría "240X80 PORCELANATO" durante el
# mes de "agosto", y propaga el efecto a las hojas resumen.

$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO": columna D (240X80 PORCELANATO), fila 23 ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D23").Value = -81.41

# --- Hoja "VENTA MENSUAL": columna F (agosto) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F23").Value = -81.41
$wsMensual.Range("F24").Value = 3137.67

# --- Hoja "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Ensanchar levemente la columna F (CUMPLIMIENTO) para que el nuevo
# porcentaje negativo quepa correctamente (de 25 a 27 caracteres).
# Nota: el setter de ColumnWidth cuantiza al ancho de píxel más cercano
# usando la métrica de la fuente; 26.17 es el valor de entrada que el
# motor resuelve exactamente a un ancho almacenado de 27.
$wsCumplimiento.Columns.Item(6).ColumnWidth = 26.17

# Fila 3: categoría "240X80 PORCELANATO"
$wsCumplimiento.Range("D3").Value = -81.41
$wsCumplimiento.Range("E3").Value = 8750.32
$wsCumplimiento.Range("F3").Value = -0.009391030706282565

# Fila 19: fila TOTAL
$wsCumplimiento.Range("D19").Value = 3137.67
$wsCumplimiento.Range("E19").Value = 51885.49386304604
$wsCumplimiento.Range("F19").Value = 0.05702452893857096
